# The "Sales / Marketing" role at "AGIS AG, Teufen" (a hidden row in the
# "Beruflicher Werdegang" / career table) was removed from the CV.
# This was stored as the hidden row 11 in the worksheet; deleting it
# shifts every following row up by one and drops its three associated
# shared strings ("Sales / Marketing", "AGIS AG, Teufen" and the bullet
# description) from the shared-strings table automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the whole row (entire row delete), shifting rows 12:16 up to 11:15.
$ws.Rows(11).Delete()

# Reflect the author's post-edit cursor position recorded in the workbook.
$ws.Range("C10").Select() | Out-Null
